$wb = $excel.ActiveWorkbook

# ALC row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 42.909092
$ws.Range("I11").Value = 42.909092
$ws.Range("K11").Value = 42.909092
$ws.Range("M11").Value = 97.090908

# ALC row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 964.45
$ws.Range("J97").Value = 1004.6842
$ws.Range("L97").Value = 3014.0526
$ws.Range("N97").Value = -4006.0526

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1418.7142
$ws.Range("I98").Value = 598
$ws.Range("J98").Value = 2896
$ws.Range("K98").Value = 598
$ws.Range("L98").Value = 2896
$ws.Range("M98").Value = 900
$ws.Range("N98").Value = -5892

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1418.7142
$ws.Range("I122").Value = 598
$ws.Range("J122").Value = 2896
$ws.Range("K122").Value = 1794
$ws.Range("L122").Value = 8688
$ws.Range("M122").Value = 656
$ws.Range("N122").Value = -13588

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1319.1072
$ws.Range("J129").Value = 1579.7727
$ws.Range("L129").Value = 4739.3181
$ws.Range("N129").Value = -14739.3181

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3144.9143
$ws.Range("I132").Value = 3582.64
$ws.Range("J132").Value = 2050.6
$ws.Range("K132").Value = 10747.92
$ws.Range("L132").Value = 6151.799999999999
$ws.Range("M132").Value = -8217.92
$ws.Range("N132").Value = -11211.8

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2286.625
$ws.Range("I141").Value = 1532.1666
$ws.Range("K141").Value = 4596.4998
$ws.Range("M141").Value = 583.5002000000004

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4819.635
$ws.Range("I32").Value = 3752.348
$ws.Range("J32").Value = 13002.167
$ws.Range("K32").Value = 3752.348
$ws.Range("L32").Value = 13002.167
$ws.Range("M32").Value = -3465.348
$ws.Range("N32").Value = -13576.167

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3281.682
$ws.Range("I45").Value = 2908.818
$ws.Range("J45").Value = 3654.5454
$ws.Range("K45").Value = 2908.818
$ws.Range("L45").Value = 3654.5454
$ws.Range("M45").Value = -2531.818
$ws.Range("N45").Value = -4408.5454

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 52633750
$ws.Range("I74").Value = 71429670
$ws.Range("K74").Value = 71429670
$ws.Range("M74").Value = -71428796

# ARM row 76
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 19998
$ws.Range("J76").Value = 19998
$ws.Range("L76").Value = 19998
$ws.Range("N76").Value = -20674

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 52633750
$ws.Range("I77").Value = 71429670
$ws.Range("K77").Value = 357148350
$ws.Range("M77").Value = -357143982

# ARM row 79
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 19998
$ws.Range("J79").Value = 19998
$ws.Range("L79").Value = 19998
$ws.Range("N79").Value = -22338

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 100001270
$ws.Range("I97").Value = 1699
$ws.Range("K97").Value = 1699
$ws.Range("M97").Value = -1203

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1109.579
$ws.Range("I99").Value = 1121.6666
$ws.Range("K99").Value = 1121.6666
$ws.Range("M99").Value = 376.3334

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4170439.2
$ws.Range("I105").Value = 5028.3335
$ws.Range("K105").Value = 5028.3335
$ws.Range("M105").Value = -3281.3335

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1028.3334
$ws.Range("I107").Value = 934
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 934
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 986
$ws.Range("N107").Value = -5340

# BSM row 118
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 43750
$ws.Range("J118").Value = 43750
$ws.Range("L118").Value = 43750
$ws.Range("N118").Value = -47064

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4310.8076
$ws.Range("I134").Value = 4403.24
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 13209.72
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -10674.72
$ws.Range("N134").Value = -11070

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2909
$ws.Range("I31").Value = 2249.7896
$ws.Range("J31").Value = 3340.8965
$ws.Range("K31").Value = 2249.7896
$ws.Range("L31").Value = 3340.8965
$ws.Range("M31").Value = -1954.7896
$ws.Range("N31").Value = -3930.8965

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2909
$ws.Range("I34").Value = 2249.7896
$ws.Range("J34").Value = 3340.8965
$ws.Range("K34").Value = 2249.7896
$ws.Range("L34").Value = 3340.8965
$ws.Range("M34").Value = -2047.7896
$ws.Range("N34").Value = -3744.8965

# CRP row 52
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 18930.482
$ws.Range("I58").Value = 1684.6666
$ws.Range("K58").Value = 1684.6666
$ws.Range("M58").Value = -1481.6666

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 33363502
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 41703130
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 41703130
$ws.Range("N86").Value = -41705376
$ws.Range("M86").Value = -3877

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 33363502
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 41703130
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 208515650
$ws.Range("N89").Value = -208526882
$ws.Range("M89").Value = -19384

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 20834586
$ws.Range("I105").Value = 41667332
$ws.Range("K105").Value = 41667332
$ws.Range("M105").Value = -41665585

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1281.1724
$ws.Range("I107").Value = 1057.2727
$ws.Range("J107").Value = 1418
$ws.Range("K107").Value = 1057.2727
$ws.Range("L107").Value = 1418
$ws.Range("M107").Value = 862.7273
$ws.Range("N107").Value = -5258

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2696.9092
$ws.Range("I132").Value = 1964.32
$ws.Range("J132").Value = 4986.25
$ws.Range("K132").Value = 5892.96
$ws.Range("L132").Value = 14958.75
$ws.Range("M132").Value = -3362.96
$ws.Range("N132").Value = -20018.75

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1140.6666
$ws.Range("I134").Value = 913.4
$ws.Range("J134").Value = 1424.75
$ws.Range("K134").Value = 2740.2
$ws.Range("L134").Value = 4274.25
$ws.Range("M134").Value = -205.1999999999998
$ws.Range("N134").Value = -9344.25

# CRP row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 50692
$ws.Range("J135").Value = 50692
$ws.Range("L135").Value = 50692
$ws.Range("N135").Value = -60832

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 18930.482
$ws.Range("I136").Value = 1684.6666
$ws.Range("K136").Value = 5053.9998
$ws.Range("M136").Value = -2503.9998

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 337.63635
$ws.Range("J23").Value = 511.42856
$ws.Range("L23").Value = 1534.28568
$ws.Range("N23").Value = -2004.28568

# CUL row 47
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 100
$ws.Range("I47").Value = 100
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 300
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 131
$ws.Range("N47").ClearContents()

# CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 3774.75
$ws.Range("I56").Value = 3774.75
$ws.Range("K56").Value = 3774.75
$ws.Range("M56").Value = -3244.75

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1091.6
$ws.Range("I117").Value = 864.5
$ws.Range("J117").Value = 2000
$ws.Range("K117").Value = 2593.5
$ws.Range("L117").Value = 6000
$ws.Range("M117").Value = 848.5
$ws.Range("N117").Value = -12884

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 269787.1
$ws.Range("I129").Value = 787.7778
$ws.Range("J129").Value = 511886.5
$ws.Range("K129").Value = 2363.3334
$ws.Range("L129").Value = 1535659.5
$ws.Range("M129").Value = 2636.6666
$ws.Range("N129").Value = -1545659.5

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 805.6836499999999
$ws.Range("I131").Value = 398.625
$ws.Range("J131").Value = 841.86664
$ws.Range("K131").Value = 1195.875
$ws.Range("L131").Value = 2525.59992
$ws.Range("M131").Value = 3844.125
$ws.Range("N131").Value = -12605.59992

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1533.6666
$ws.Range("I97").Value = 1665
$ws.Range("K97").Value = 1665
$ws.Range("M97").Value = -1169

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 25003098
$ws.Range("I102").Value = 27780580
$ws.Range("K102").Value = 27780580
$ws.Range("M102").Value = -27778958

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3262.6667
$ws.Range("J40").Value = 5270
$ws.Range("L40").Value = 5270
$ws.Range("N40").Value = -5542

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 641.05
$ws.Range("I46").Value = 559.1818
$ws.Range("J46").Value = 741.1111
$ws.Range("K46").Value = 559.1818
$ws.Range("L46").Value = 741.1111
$ws.Range("M46").Value = -371.1818
$ws.Range("N46").Value = -1117.1111

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2108.5557
$ws.Range("I100").Value = 996.6667
$ws.Range("J100").Value = 2664.5
$ws.Range("K100").Value = 996.6667
$ws.Range("L100").Value = 2664.5
$ws.Range("M100").Value = -455.6667
$ws.Range("N100").Value = -3746.5

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 66667710
$ws.Range("I81").Value = 1166.3
$ws.Range("K81").Value = 2332.6
$ws.Range("M81").Value = -1271.6

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 66667710
$ws.Range("I84").Value = 1166.3
$ws.Range("K84").Value = 11663
$ws.Range("M84").Value = -6359

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 270
$ws.Range("I100").Value = 270
$ws.Range("K100").Value = 540
$ws.Range("M100").Value = 1

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2458703
$ws.Range("I113").Value = 2047.3334
$ws.Range("J113").Value = 13513654
$ws.Range("K113").Value = 6142.0002
$ws.Range("L113").Value = 40540962
$ws.Range("M113").Value = -3972.0002
$ws.Range("N113").Value = -40545302

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 22940538
$ws.Range("J136").Value = 2800.5557
$ws.Range("L136").Value = 8401.667099999999
$ws.Range("N136").Value = -13501.6671
